$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force Text format on price cells whose new values would otherwise be
# auto-coerced to numbers by Excel (losing the literal-text representation).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "47.674.15"
$ws.Range("E2").Value = "  +5.24%  "

$ws.Range("D3").Value = "2.640.45"
$ws.Range("E3").Value = "  +8.78%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "312.01"
$ws.Range("E5").Value = "  +6.29%  "

$ws.Range("D6").Value = "103.35"
$ws.Range("E6").Value = "  +10.17%  "

$ws.Range("D7").Value = "0.609"
$ws.Range("E7").Value = "  +9.15%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "0.596"
$ws.Range("E9").Value = "  +18.99%  "

$ws.Range("D10").Value = "40.34"
$ws.Range("E10").Value = "  +18.06%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.0862"
$ws.Range("E11").Value = "  +10.83%  "

$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").Value = "55.41"
$ws.Range("E12").Value = "  +3.28%  "

$ws.Range("D13").Value = "8.33"
$ws.Range("E13").Value = "  +18.96%  "

$ws.Range("D14").Value = "3.053.63"
$ws.Range("E14").Value = "  +9.11%  "

$ws.Range("E15").Value = "  +2.69%  "

$ws.Range("D16").Value = "2.656.15"
$ws.Range("E16").Value = "  +8.92%  "

$ws.Range("E17").Value = "  +11.93%  "

$ws.Range("E18").Value = "  +7.56%  "

$ws.Range("D19").Value = "47.778.46"
$ws.Range("E19").Value = "  +5.59%  "

$ws.Range("E20").Value = "  +10.18%  "

$ws.Range("D21").Value = "13.24"
$ws.Range("E21").Value = "  +6.97%  "

$ws.Range("D22").Value = "6.83"
$ws.Range("E22").Value = "  +10.37%  "

$ws.Range("D23").Value = "73.41"
$ws.Range("E23").Value = "  +9.45%  "

$ws.Range("D24").Value = "280.76"
$ws.Range("E24").Value = "  +17.35%  "

$ws.Range("E25").Value = "  +11.58%  "

$ws.Range("D26").Value = "30.88"
$ws.Range("E26").Value = "  +44.73%  "

$ws.Range("E27").Value = "  +15.42%  "

$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.18%  "

$ws.Range("E29").Value = "  +1.40%  "

$ws.Range("D30").Value = "10.71"
$ws.Range("E30").Value = "  +12.22%  "

$ws.Range("D31").Value = "2.34"
$ws.Range("E31").Value = "  +5.22%  "

$ws.Range("D32").Value = "39.87"
$ws.Range("E32").Value = "  +7.21%  "

$ws.Range("D33").Value = "6.24"
$ws.Range("E33").Value = "  +15.50%  "

$ws.Range("D34").Value = "3.71"
$ws.Range("E34").Value = "  -3.90%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.0858"
$ws.Range("E35").Value = "  +12.79%  "

$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "2.87"
$ws.Range("E36").Value = "  +5.69%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.129"
$ws.Range("E37").Value = "  +14.74%  "

$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "2.25"
$ws.Range("E38").Value = "  +12.62%  "

$ws.Range("D39").Value = "153.16"
$ws.Range("E39").Value = "  +2.67%  "

$ws.Range("D40").Value = "0.125"
$ws.Range("E40").Value = "  +8.97%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "4.26"
$ws.Range("E41").Value = "  +14.53%  "

$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "16.09"
$ws.Range("E42").Value = "  +12.54%  "

$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").Value = "3.73"
$ws.Range("E43").Value = "  +17.87%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "22.58"
$ws.Range("E44").Value = "  +38.82%  "

$ws.Range("D45").Value = "0.0336"
$ws.Range("E45").Value = "  +14.49%  "

$ws.Range("D46").Value = "2.176.52"
$ws.Range("E46").Value = "  +9.16%  "

$ws.Range("D47").Value = "97.92"
$ws.Range("E47").Value = "  +10.80%  "

$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.10%  "

$ws.Range("D49").Value = "9.85"
$ws.Range("E49").Value = "  +14.93%  "

$ws.Range("D50").Value = "1.84"
$ws.Range("E50").Value = "  +8.56%  "

$ws.Range("D51").Value = "110.26"
$ws.Range("E51").Value = "  +9.00%  "
